# Add a new row (row 15) to Sheet1 containing a newly documented D3D12 error
# ("Textures created with certain Formats must align the resource dimensions
# properly...") and its Korean-language resolution note, mirroring the other
# error/solution rows already present on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$errorText = "Textures created with certain Formats must align the resource dimensions properly. D3D12_RESOURCE_DESC::Format is BC3_UNORM. D3D12_RESOURCE_DESC::Width is 512, and must be a multiple of 4. D3D12_RESOURCE_DESC::Height is 170, and must be a multiple of 4. [ STATE_CREATION ERROR #597: CREATERESOURCE_INVALIDDIMENSIONS]"
$solutionText = "밉맵추출을 위해 dds는 4의 배수 텍스쳐야 한다."

$ws.Range("A15").Value2 = $errorText
$ws.Range("B15").Value2 = $solutionText

# Match the row height used by the other wrapped-text rows on the sheet (33pt).
$ws.Rows.Item(15).RowHeight = 33

# Keep the existing selection on B15 (matches the sheet's prior selection).
$ws.Range("B15").Select() | Out-Null
